$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Grade Earned inputs
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.78
$ws.Range("G10").Value = 0.6
$ws.Range("D15").Value = 0.9

# Update the active selection to match the new cursor position
$ws.Range("G10").Select()
